$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A180").Value = "2023-12-11 11:02:23"
$ws.Range("B180").Value = 0.003

$ws.Range("A181").Value = "2023-12-11 11:02:41"
$ws.Range("B181").Value = 0.001
